$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 8).Value = 4000  # H4
$ws.Cells.Item(4, 9).Value = 4000  # I4
$ws.Cells.Item(4, 11).Value = 4000  # K4
$ws.Cells.Item(4, 13).Value = -3886  # M4
$ws.Cells.Item(17, 8).Value = 870713.1  # H17
$ws.Cells.Item(17, 10).Value = 1044760.75  # J17
$ws.Cells.Item(17, 12).Value = 3134282.25  # L17
$ws.Cells.Item(17, 14).Value = -3134618.25  # N17
$ws.Cells.Item(43, 8).Value = 10983.833  # H43
$ws.Cells.Item(43, 9).Value = 4000  # I43
$ws.Cells.Item(43, 10).Value = 11618.728  # J43
$ws.Cells.Item(43, 11).Value = 4000  # K43
$ws.Cells.Item(43, 12).Value = 11618.728  # L43
$ws.Cells.Item(43, 13).Value = -3931  # M43
$ws.Cells.Item(43, 14).Value = -11756.728  # N43
$ws.Cells.Item(135, 8).Value = 654.15  # H135
$ws.Cells.Item(135, 9).Value = 635.9474  # I135
$ws.Cells.Item(135, 11).Value = 5723.5266  # K135
$ws.Cells.Item(135, 13).Value = -3188.5266  # M135

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2243.6545  # H32
$ws.Cells.Item(32, 9).Value = 1567.94  # I32
$ws.Cells.Item(32, 11).Value = 1567.94  # K32
$ws.Cells.Item(32, 13).Value = -1280.94  # M32
$ws.Cells.Item(45, 8).Value = 55560110  # H45
$ws.Cells.Item(45, 9).Value = 76925860  # I45
$ws.Cells.Item(45, 11).Value = 76925860  # K45
$ws.Cells.Item(45, 13).Value = -76925483  # M45
$ws.Cells.Item(61, 8).Value = 7126.3125  # H61
$ws.Cells.Item(61, 9).Value = 6097.6313  # I61
$ws.Cells.Item(61, 11).Value = 6097.6313  # K61
$ws.Cells.Item(61, 13).Value = -5885.6313  # M61
$ws.Cells.Item(122, 8).Value = 4597.391  # H122
$ws.Cells.Item(122, 9).Value = 3702.2144  # I122
$ws.Cells.Item(122, 11).Value = 11106.6432  # K122
$ws.Cells.Item(122, 13).Value = -8656.643199999999  # M122
$ws.Cells.Item(136, 8).Value = 7126.3125  # H136
$ws.Cells.Item(136, 9).Value = 6097.6313  # I136
$ws.Cells.Item(136, 11).Value = 18292.8939  # K136
$ws.Cells.Item(136, 13).Value = -15742.8939  # M136

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 3877.457  # H134
$ws.Cells.Item(134, 9).Value = 1987.7693  # I134
$ws.Cells.Item(134, 11).Value = 5963.3079  # K134
$ws.Cells.Item(134, 13).Value = -3428.3079  # M134

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 22579.855  # H31
$ws.Cells.Item(31, 9).Value = 3410.7058  # I31
$ws.Cells.Item(31, 11).Value = 3410.7058  # K31
$ws.Cells.Item(31, 13).Value = -3115.7058  # M31
$ws.Cells.Item(34, 8).Value = 22579.855  # H34
$ws.Cells.Item(34, 9).Value = 3410.7058  # I34
$ws.Cells.Item(34, 11).Value = 3410.7058  # K34
$ws.Cells.Item(34, 13).Value = -3208.7058  # M34
$ws.Cells.Item(86, 8).Value = 5249.2256  # H86
$ws.Cells.Item(86, 9).Value = 3577.2144  # I86
$ws.Cells.Item(86, 10).Value = 6626.1763  # J86
$ws.Cells.Item(86, 11).Value = 3577.2144  # K86
$ws.Cells.Item(86, 12).Value = 6626.1763  # L86
$ws.Cells.Item(86, 13).Value = -2454.2144  # M86
$ws.Cells.Item(86, 14).Value = -8872.176299999999  # N86
$ws.Cells.Item(89, 8).Value = 5249.2256  # H89
$ws.Cells.Item(89, 9).Value = 3577.2144  # I89
$ws.Cells.Item(89, 10).Value = 6626.1763  # J89
$ws.Cells.Item(89, 11).Value = 17886.072  # K89
$ws.Cells.Item(89, 12).Value = 33130.8815  # L89
$ws.Cells.Item(89, 13).Value = -12270.072  # M89
$ws.Cells.Item(89, 14).Value = -44362.8815  # N89
$ws.Cells.Item(94, 8).Value = 4822.5  # H94
$ws.Cells.Item(94, 10).Value = 4096.6665  # J94
$ws.Cells.Item(94, 12).Value = 4096.6665  # L94
$ws.Cells.Item(94, 14).Value = -4998.6665  # N94
$ws.Cells.Item(122, 8).Value = 7269.3213  # H122
$ws.Cells.Item(122, 9).Value = 4176.2354  # I122
$ws.Cells.Item(122, 11).Value = 12528.7062  # K122
$ws.Cells.Item(122, 13).Value = -10078.7062  # M122
$ws.Cells.Item(132, 8).Value = 3983.739  # H132
$ws.Cells.Item(132, 9).Value = 3338.875  # I132
$ws.Cells.Item(132, 11).Value = 10016.625  # K132
$ws.Cells.Item(132, 13).Value = -7486.625  # M132

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 3086.923  # H68
$ws.Cells.Item(68, 10).Value = 3150.4  # J68
$ws.Cells.Item(68, 12).Value = 9451.200000000001  # L68
$ws.Cells.Item(68, 14).Value = -11073.2  # N68
$ws.Cells.Item(71, 8).Value = 3086.923  # H71
$ws.Cells.Item(71, 10).Value = 3150.4  # J71
$ws.Cells.Item(71, 12).Value = 28353.6  # L71
$ws.Cells.Item(71, 14).Value = -36465.60000000001  # N71
$ws.Cells.Item(113, 8).Value = 1539.8889  # H113
$ws.Cells.Item(113, 10).Value = 1919.8  # J113
$ws.Cells.Item(113, 12).Value = 5759.4  # L113
$ws.Cells.Item(113, 14).Value = -10099.4  # N113
$ws.Cells.Item(132, 8).Value = 3811.8965  # H132
$ws.Cells.Item(132, 9).Value = 2703.5  # I132
$ws.Cells.Item(132, 11).Value = 24331.5  # K132
$ws.Cells.Item(132, 13).Value = -21801.5  # M132

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(46, 8).Value = 1272.6364  # H46
$ws.Cells.Item(46, 9).Value = 1272.6364  # I46
$ws.Cells.Item(46, 11).Value = 1272.6364  # K46
$ws.Cells.Item(46, 13).Value = -1116.6364  # M46
$ws.Cells.Item(102, 8).Value = 3765.5833  # H102
$ws.Cells.Item(102, 9).Value = 3031.1667  # I102
$ws.Cells.Item(102, 11).Value = 3031.1667  # K102
$ws.Cells.Item(102, 13).Value = -1409.1667  # M102
$ws.Cells.Item(122, 8).Value = 8763.134  # H122
$ws.Cells.Item(122, 9).Value = 4493.5  # I122
$ws.Cells.Item(122, 11).Value = 13480.5  # K122
$ws.Cells.Item(122, 13).Value = -11030.5  # M122

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 3812.25  # H16
$ws.Cells.Item(16, 9).Value = 5500.5  # I16
$ws.Cells.Item(16, 11).Value = 5500.5  # K16
$ws.Cells.Item(16, 13).Value = -5330.5  # M16
$ws.Cells.Item(20, 8).Value = 19033.334  # H20
$ws.Cells.Item(20, 9).Value = 19033.334  # I20
$ws.Cells.Item(20, 11).Value = 19033.334  # K20
$ws.Cells.Item(20, 13).Value = -18807.334  # M20
$ws.Cells.Item(25, 8).Value = 5000  # H25
$ws.Cells.Item(25, 9).Value = 5000  # I25
$ws.Cells.Item(25, 11).Value = 5000  # K25
$ws.Cells.Item(25, 13).Value = -4770  # M25
$ws.Cells.Item(40, 8).Value = 7661.2354  # H40
$ws.Cells.Item(40, 9).Value = 5294.364  # I40
$ws.Cells.Item(40, 11).Value = 5294.364  # K40
$ws.Cells.Item(40, 13).Value = -5158.364  # M40
$ws.Cells.Item(63, 8).Value = 45000  # H63
$ws.Cells.Item(63, 9).Value = 45000  # I63
$ws.Cells.Item(63, 11).Value = 45000  # K63
$ws.Cells.Item(63, 13).Value = -44251  # M63
$ws.Cells.Item(66, 8).Value = 45000  # H66
$ws.Cells.Item(66, 9).Value = 45000  # I66
$ws.Cells.Item(66, 11).Value = 135000  # K66
$ws.Cells.Item(66, 13).Value = -131256  # M66
$ws.Cells.Item(82, 8).Value = 21908.4  # H82
$ws.Cells.Item(82, 9).Value = 18997.334  # I82
$ws.Cells.Item(82, 10).Value = 26275  # J82
$ws.Cells.Item(82, 11).Value = 18997.334  # K82
$ws.Cells.Item(82, 12).Value = 26275  # L82
$ws.Cells.Item(82, 13).Value = -18636.334  # M82
$ws.Cells.Item(82, 14).Value = -26997  # N82
$ws.Cells.Item(85, 8).Value = 21908.4  # H85
$ws.Cells.Item(85, 9).Value = 18997.334  # I85
$ws.Cells.Item(85, 10).Value = 26275  # J85
$ws.Cells.Item(85, 11).Value = 18997.334  # K85
$ws.Cells.Item(85, 12).Value = 26275  # L85
$ws.Cells.Item(85, 13).Value = -17749.334  # M85
$ws.Cells.Item(85, 14).Value = -28771  # N85
$ws.Cells.Item(93, 8).Value = 43422.582  # H93
$ws.Cells.Item(93, 9).Value = 31564.834  # I93
$ws.Cells.Item(93, 11).Value = 31564.834  # K93
$ws.Cells.Item(93, 13).Value = -30316.834  # M93
$ws.Cells.Item(122, 8).Value = 7244.9473  # H122
$ws.Cells.Item(122, 9).Value = 5874.643  # I122
$ws.Cells.Item(122, 11).Value = 17623.929  # K122
$ws.Cells.Item(122, 13).Value = -15173.929  # M122
$ws.Cells.Item(136, 8).Value = 5486.6665  # H136
$ws.Cells.Item(136, 9).Value = 6642.5557  # I136
$ws.Cells.Item(136, 11).Value = 19927.6671  # K136
$ws.Cells.Item(136, 13).Value = -17377.6671  # M136

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(18, 8).Value = 13907  # H18
$ws.Cells.Item(18, 9).Value = 0  # I18
$ws.Cells.Item(18, 10).Value = 13907  # J18
$ws.Cells.Item(18, 11).Value = 0  # K18
$ws.Cells.Item(18, 12).Value = 13907  # L18
$ws.Cells.Item(18, 13).ClearContents()  # M18
$ws.Cells.Item(18, 14).Value = -14253  # N18
$ws.Cells.Item(52, 8).Value = 12000  # H52
$ws.Cells.Item(52, 10).Value = 11250  # J52
$ws.Cells.Item(52, 12).Value = 11250  # L52
$ws.Cells.Item(52, 14).Value = -11702  # N52
$ws.Cells.Item(122, 8).Value = 4185.1665  # H122
$ws.Cells.Item(122, 9).Value = 1859.6666  # I122
$ws.Cells.Item(122, 10).Value = 9611.333000000001  # J122
$ws.Cells.Item(122, 11).Value = 5578.9998  # K122
$ws.Cells.Item(122, 12).Value = 28833.999  # L122
$ws.Cells.Item(122, 13).Value = -3128.9998  # M122
$ws.Cells.Item(122, 14).Value = -33733.999  # N122
$ws.Cells.Item(136, 8).Value = 3967.9473  # H136
$ws.Cells.Item(136, 9).Value = 2212.9666  # I136
$ws.Cells.Item(136, 11).Value = 6638.899800000001  # K136
$ws.Cells.Item(136, 13).Value = -4088.899800000001  # M136
$ws.Cells.Item(139, 8).Value = 69443.55499999999  # H139
$ws.Cells.Item(139, 10).Value = 69443.55499999999  # J139
$ws.Cells.Item(139, 12).Value = 69443.55499999999  # L139
$ws.Cells.Item(139, 14).Value = -79723.55499999999  # N139
$ws.Cells.Item(141, 8).Value = 275715  # H141
$ws.Cells.Item(141, 10).Value = 275715  # J141
$ws.Cells.Item(141, 12).Value = 275715  # L141
$ws.Cells.Item(141, 14).Value = -286075  # N141
